try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # Insert a new row at position 29 (shifts existing rows 29-93 down to 30-94,
    # preserving all their data/formatting as-is).
    $ws.Rows(29).Insert()

    # Populate the newly inserted row 29 with the new data entry.
    $ws.Cells.Item(29, 1).Value = 5
    $ws.Cells.Item(29, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item(29, 3).Value = "Maule"
    $ws.Cells.Item(29, 4).Value = 44526
    $ws.Cells.Item(29, 5).Value = 7
    $ws.Cells.Item(29, 6).Value = "Fruta"
    $ws.Cells.Item(29, 7).Value = 100108
    $ws.Cells.Item(29, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item(29, 9).Value = 100108002
    $ws.Cells.Item(29, 10).Value = "Mango"
    $ws.Cells.Item(29, 11).Value = "Sin especificar"
    $ws.Cells.Item(29, 12).Value = "Primera"
    $ws.Cells.Item(29, 13).Value = 300
    $ws.Cells.Item(29, 14).Value = 6000
    $ws.Cells.Item(29, 15).Value = 6000
    $ws.Cells.Item(29, 16).Value = 6000
    $ws.Cells.Item(29, 17).Value = "`$/bandeja 4 kilos"
    $ws.Cells.Item(29, 18).Value = "Perú"
    $ws.Cells.Item(29, 19).Value = 1500
    $ws.Cells.Item(29, 20).Value = 4
} catch {
    Write-Host "ERROR: $_"
}
